$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo: "Fling" -> "Flint" in the product name (cell B4)
$ws.Range("B4").Value = "Flint Match Lighter Metal Outdoor Camping"

# Update the active selection to match the saved view state (C4 -> B4)
$ws.Range("B4").Select()
